# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.228.06'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.860.18'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7157'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07758'
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.18'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08257'
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.238'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.847.14'
$ws.Range("E13").Value = '  -1.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7177'
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.28'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.235.85'
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '244.28'
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007797'
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.108.79'
$ws.Range("E21").Value = '  -2.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.973'
$ws.Range("E23").Value = '  +3.02%  '
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1602'
$ws.Range("E25").Value = '  +2.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.50'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.497'
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("E30").Value = '  -3.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.407'
$ws.Range("E31").Value = '  +1.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.188'
$ws.Range("E32").Value = '  +2.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05194'
$ws.Range("E33").Value = '  -1.04%  '
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("E35").Value = '  -2.03%  '
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.678'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01857'
$ws.Range("E38").Value = '  -0.65%  '
$ws.Range("E39").Value = '  -1.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.148.28'
$ws.Range("E40").Value = '  -2.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9064'
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.128'
$ws.Range("E42").Value = '  +2.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.27'
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9995'
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.81'
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.004.68'
$ws.Range("E46").Value = '  -2.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5230'
$ws.Range("E47").Value = '  -2.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.771'
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.327'
$ws.Range("E50").Value = '  +1.79%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.875'
$ws.Range("E51").Value = '  +1.41%  '
